$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).Value = 37
$ws.Cells.Item($row, 2).Value = "croatia"
$ws.Cells.Item($row, 3).Value = "hnl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45191.75
$ws.Cells.Item($row, 6).Value = "Osijek"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Varazdin"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 1.54
$ws.Cells.Item($row, 11).Value = "17/09/2023 18:42"
$ws.Cells.Item($row, 12).Value = 1.65
$ws.Cells.Item($row, 13).Value = "22/09/2023 17:51"
$ws.Cells.Item($row, 14).Value = 4.25
$ws.Cells.Item($row, 15).Value = "17/09/2023 18:42"
$ws.Cells.Item($row, 16).Value = 4.13
$ws.Cells.Item($row, 17).Value = "22/09/2023 17:51"
$ws.Cells.Item($row, 18).Value = 5.15
$ws.Cells.Item($row, 19).Value = "17/09/2023 18:42"
$ws.Cells.Item($row, 20).Value = 4.97
$ws.Cells.Item($row, 21).Value = "22/09/2023 17:51"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/osijek-varazdin/lvKZhpie/"

# Match formatting of the previous row (A: bold/border/centered index cell,
# E: date-time number format) by copying formats from row 37 onto row 38.
$ws.Range("A37").Copy()
$ws.Range("A38").PasteSpecial(-4122)

$ws.Range("E37").Copy()
$ws.Range("E38").PasteSpecial(-4122)

$excel.CutCopyMode = 0
